$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'27.905.53"
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.Value = "'  +0.81%  "
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.Value = "'1.880.00"
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.Value = "'  +0.05%  "
$c.Style = "Normal"
$c = $ws.Range("E4")
$c.Value = "'  +1.73%  "
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.Value = "'334.49"
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.Value = "'  +0.99%  "
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.Value = "'1.018"
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.Value = "'  +1.57%  "
$c.Style = "Normal"
$c = $ws.Range("E7")
$c.Value = "'  -1.71%  "
$c.Style = "Normal"
$c = $ws.Range("D8")
$c.Value = "'0.3896"
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.Value = "'  -1.85%  "
$c.Style = "Normal"
$c = $ws.Range("D9")
$c.Value = "'46.92"
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.Value = "'  -1.85%  "
$c.Style = "Normal"
$c = $ws.Range("D10")
$c.Value = "'0.07925"
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.Value = "'  -1.66%  "
$c.Style = "Normal"
$c = $ws.Range("D11")
$c.Value = "'1.005"
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.Value = "'  -1.86%  "
$c.Style = "Normal"
$c = $ws.Range("D12")
$c.Value = "'21.53"
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.Value = "'  -1.47%  "
$c.Style = "Normal"
$c = $ws.Range("D13")
$c.Value = "'1.901.79"
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.Value = "'  +2.54%  "
$c.Style = "Normal"
$c = $ws.Range("D14")
$c.Value = "'5.913"
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.Value = "'  -0.90%  "
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.Value = "'7.060"
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.Value = "'  -1.52%  "
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.Value = "'  +1.90%  "
$c.Style = "Normal"
$c = $ws.Range("D17")
$c.Value = "'0.06758"
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.Value = "'  +2.06%  "
$c.Style = "Normal"
$c = $ws.Range("D18")
$c.Value = "'86.71"
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.Value = "'  -0.50%  "
$c.Style = "Normal"
$c = $ws.Range("D19")
$c.Value = "'0.00001039"
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.Value = "'  -0.28%  "
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.Value = "'  -1.40%  "
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.Value = "'  +1.64%  "
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.Value = "'27.917.02"
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.Value = "'  +0.81%  "
$c.Style = "Normal"
$c = $ws.Range("D23")
$c.Value = "'5.450"
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.Value = "'  -1.11%  "
$c.Style = "Normal"
$c = $ws.Range("D24")
$c.Value = "'10.88"
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.Value = "'  -1.36%  "
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.Value = "'  +2.30%  "
$c.Style = "Normal"
$c = $ws.Range("D26")
$c.Value = "'2.119.04"
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.Value = "'  +1.90%  "
$c.Style = "Normal"
$c = $ws.Range("D27")
$c.Value = "'159.59"
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.Value = "'  +2.01%  "
$c.Style = "Normal"
$c = $ws.Range("D28")
$c.Value = "'19.84"
$c.Style = "Normal"
$c = $ws.Range("E28")
$c.Value = "'  -2.13%  "
$c.Style = "Normal"
$c = $ws.Range("E29")
$c.Value = "'  -2.16%  "
$c.Style = "Normal"
$c = $ws.Range("D30")
$c.Value = "'5.423"
$c.Style = "Normal"
$c = $ws.Range("E30")
$c.Value = "'  -3.18%  "
$c.Style = "Normal"
$c = $ws.Range("D31")
$c.Value = "'120.82"
$c.Style = "Normal"
$c = $ws.Range("E31")
$c.Value = "'  -1.32%  "
$c.Style = "Normal"
$c = $ws.Range("D32")
$c.Value = "'0.09458"
$c.Style = "Normal"
$c = $ws.Range("E32")
$c.Value = "'  -1.19%  "
$c.Style = "Normal"
$c = $ws.Range("D33")
$c.Value = "'0.9534"
$c.Style = "Normal"
$c = $ws.Range("E33")
$c.Value = "'  -1.99%  "
$c.Style = "Normal"
$c = $ws.Range("D34")
$c.Value = "'3.665"
$c.Style = "Normal"
$c = $ws.Range("E34")
$c.Value = "'  +1.09%  "
$c.Style = "Normal"
$c = $ws.Range("D35")
$c.Value = "'5.290"
$c.Style = "Normal"
$c = $ws.Range("E35")
$c.Value = "'  -0.56%  "
$c.Style = "Normal"
$c = $ws.Range("E36")
$c.Value = "'  -7.09%  "
$c.Style = "Normal"
$c = $ws.Range("D37")
$c.Value = "'0.06095"
$c.Style = "Normal"
$c = $ws.Range("E37")
$c.Value = "'  -0.43%  "
$c.Style = "Normal"
$c = $ws.Range("D38")
$c.Value = "'0.02226"
$c.Style = "Normal"
$c = $ws.Range("E38")
$c.Value = "'  -1.55%  "
$c.Style = "Normal"
$c = $ws.Range("D39")
$c.Value = "'1.206"
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.Value = "'  -2.14%  "
$c.Style = "Normal"
$c = $ws.Range("D40")
$c.Value = "'8.064"
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.Value = "'  -1.13%  "
$c.Style = "Normal"
$c = $ws.Range("D41")
$c.Value = "'0.5860"
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.Value = "'  -2.68%  "
$c.Style = "Normal"
$c = $ws.Range("D42")
$c.Value = "'0.1878"
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.Value = "'  -1.41%  "
$c.Style = "Normal"
$c = $ws.Range("D43")
$c.Value = "'10.08"
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.Value = "'  -1.73%  "
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.Value = "'  +1.20%  "
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.Value = "'0.5605"
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.Value = "'  -1.76%  "
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.Value = "'  -1.56%  "
$c.Style = "Normal"
$c = $ws.Range("D47")
$c.Value = "'3.379"
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.Value = "'  -0.80%  "
$c.Style = "Normal"
$c = $ws.Range("D48")
$c.Value = "'1.903"
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.Value = "'  -1.57%  "
$c.Style = "Normal"
$c = $ws.Range("D49")
$c.Value = "'0.06896"
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.Value = "'  +1.10%  "
$c.Style = "Normal"
$c = $ws.Range("D50")
$c.Value = "'113.21"
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.Value = "'  +1.64%  "
$c.Style = "Normal"
$c = $ws.Range("B51")
$c.Value = "'EOS"
$c.Style = "Normal"
$c = $ws.Range("C51")
$c.Value = "'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$c.Style = "Normal"
$c = $ws.Range("D51")
$c.Value = "'1.058"
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.Value = "'  -1.43%  "
$c.Style = "Normal"
